$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("L2").Value = "[4.264382927105267, 8.77939936043948]"
$ws.Range("M2").Value = [double]"2.672089505217912e-08"
$ws.Range("N2").Value = [double]"5.344179010435823e-08"
$ws.Range("P2").Value = "[-1.8491055858966945, -1.0440528138056173]"
$ws.Range("Q2").Value = [double]"7.635003740347202e-12"
$ws.Range("R2").Value = [double]"7.635003740347202e-12"
$ws.Range("T2").Value = "[7.968339124338804, 10.727283164703103]"
$ws.Range("X2").Value = 4.238898898898992
$ws.Range("Y2").Value = 7.507447447447607

# Row 3 updates
$ws.Range("L3").Value = "[4.204186229752095, 9.199380381080681]"
$ws.Range("M3").Value = [double]"2.545926187647041e-07"
$ws.Range("N3").Value = [double]"2.545926187647041e-07"
$ws.Range("P3").Value = "[-3.434053230951005, -2.528368862348542]"
$ws.Range("T3").Value = "[7.101728817330797, 10.098647936492252]"
$ws.Range("X3").Value = 9.899099099099262
$ws.Range("Y3").Value = 13.44504504504527
